# Commit: "Changed hospital to health facility"
#
# This script:
#  1. Renames the "...hospital..." labels to "...health facility..." everywhere
#     in the workbook ("IFAS (hospital)" -> "IFAS (health facility)" and
#     "IFAS for pregnant women (hospital)" -> "IFAS for pregnant women (health facility)").
#  2. Relabels the unit-cost column header on "Programs cost and coverage"
#     from "Unit cost (US$)" to "Unit cost (US$ per person per year)".
#  3. Updates two baseline staple-food fractions (rice / wheat) from 0 to 0.1
#     on "Baseline year population inputs" (dependent formulas recalculate
#     automatically).
#  4. Converts two formula-driven unit-cost cells on "Programs cost and
#     coverage" (diarrhoea treatment D6 and SAM treatment D30) into plain
#     hard-coded values, matching the formatting of the other hard-coded
#     unit-cost cells in that column.
#  5. Adds Nick Scott's explanatory cell comments to D28, D30 and D37 on
#     "Programs cost and coverage".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "hospital" -> "health facility" label renames (applies to every sheet
#    that references these shared strings: "Programs cost and coverage",
#    "Programs target population", "Program dependencies", etc.)
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("IFAS for pregnant women (hospital)", "IFAS for pregnant women (health facility)") | Out-Null
    $ws.Cells.Replace("IFAS (hospital)", "IFAS (health facility)") | Out-Null
}

$progCost = $wb.Worksheets.Item("Programs cost and coverage")

# ---------------------------------------------------------------------------
# 2. Unit cost column header relabel
# ---------------------------------------------------------------------------
$progCost.Range("D1").Value = "Unit cost (US`$ per person per year)"

# ---------------------------------------------------------------------------
# 3. Baseline staple food fractions
# ---------------------------------------------------------------------------
$baseline = $wb.Worksheets.Item("Baseline year population inputs")
$baseline.Range("C16").Value = 0.1
$baseline.Range("C17").Value = 0.1

# ---------------------------------------------------------------------------
# 4. Replace linked-formula unit costs with hard-coded values (and restore
#    the plain "input" formatting used by the other unit-cost cells, e.g.
#    D2, instead of the "calculated" formatting the formula cells carried).
# ---------------------------------------------------------------------------
$progCost.Range("D2").Copy() | Out-Null
$progCost.Range("D6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$progCost.Range("D6").Value = 0.82

$progCost.Range("D2").Copy() | Out-Null
$progCost.Range("D30").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$progCost.Range("D30").Value = 5.3

# ---------------------------------------------------------------------------
# 5. Nick Scott's comments
# ---------------------------------------------------------------------------
$commentD28 = "Nick Scott:`nThe cost per child per year can be estimated as `n= (cost per treatment) * (annual diarrhoea incidence)`n`nDiarrhoea incidence is the average in children under 5.  See user guide for further information"
$progCost.Range("D28").AddComment($commentD28) | Out-Null

$commentD30 = "Nick Scott:`nThe cost per child per year can be estimated as `n= (cost per treatment episode) * (SAM prevalence) * 2.6`nCost per treatment episode includes management of MAM (if selected) and is an average over delivery modalities. See user guide for further information"
$progCost.Range("D30").AddComment($commentD30) | Out-Null

$commentD37 = "Nick Scott:`nThe cost per child per year can be estimated as `n= (cost per treatment) * (annual diarrhoea incidence)`nDiarrhoea incidence is the average in children under 5.  See user guide for further information"
$progCost.Range("D37").AddComment($commentD37) | Out-Null

Write-Output "edit complete"
